$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# The old rows 11 ("up") and 12 ("dn") are removed entirely; this shifts every
# row below them up by two (rows 13-22 become rows 11-20), preserving all
# existing per-row formatting (e.g. the multi-line "Gain/Phase/Auto" rows).
$ws.Rows("11:12").Delete()

# Text already in place for D4/F4/F5/F8/F15/F18 needs updating for the new
# single-value ("digit") frequency/phase-increment entry widgets.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "Center"
$ws.Range("D4").Style = "Normal"

$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "<FVAL>"
$ws.Range("F4").Style = "Normal"

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "Frequency  [Hertz]"
$ws.Range("F5").Style = "Normal"

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "phase inc: <value>"
$ws.Range("F8").Style = "Normal"

$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "9876543.210"
$ws.Range("F15").Style = "Normal"

$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "FFEEDDCC"
$ws.Range("F18").Style = "Normal"

# Brand new row 21 ("Back" button text id).
$ws.Range("B21:F21").NumberFormat = "@"
$ws.Range("B21").Value = "SingleUseId40"
$ws.Range("C21").Value = "Default"
$ws.Range("D21").Value = "Left"
$ws.Range("E21").Value = "LTR"
$ws.Range("F21").Value = "Back"
$ws.Range("B21:F21").Style = "Normal"
